$p = $ppt.ActivePresentation

# --- 1) Table on slide 5: switch its table style to the built-in
#        "Light Style 2" table style ({AB1C0CA1-0F41-4ADB-89C6-8ECFC1D4ADA1}).
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{AB1C0CA1-0F41-4ADB-89C6-8ECFC1D4ADA1}")

# --- 2) Re-colour the deck's theme (the "Integral"/"Red Violet" scheme)
#        to the standard Office theme colours.
$themeColors = $p.Slides.Item(1).Master.Theme.ThemeColorScheme
$themeColors.Item(1).RGB  = 0         # dk1      000000
$themeColors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388   # dk2      44546A
$themeColors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407     # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Item(10).RGB = 4697456   # accent6  70AD47
$themeColors.Item(11).RGB = 12673797  # hlink    0563C1
$themeColors.Item(12).RGB = 7491477   # folHlink 954F72
